$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row updates: column letter -> new text value.
# D-column values are forced to Text format before assignment so Excel
# doesn't auto-convert numeric-looking strings (e.g. "251.85") into numbers,
# preserving the original inline-string text type of these cells.
$updates = @{
    2  = @{ D = "30.350.28"; E = "  +0.33%  " }
    3  = @{ D = "1.938.02"; E = "  +0.59%  " }
    5  = @{ D = "251.85" }
    6  = @{ D = "0.7243"; E = "  +3.46%  " }
    7  = @{ D = "1.000"; E = "  -0.16%  " }
    8  = @{ D = "0.3319"; E = "  +2.67%  " }
    9  = @{ D = "27.98"; E = "  +5.87%  " }
    10 = @{ D = "0.07273"; E = "  +6.66%  " }
    11 = @{ E = "  +1.96%  " }
    12 = @{ D = "0.08108"; E = "  +2.10%  " }
    13 = @{ D = "1.934.21"; E = "  +0.43%  " }
    14 = @{ D = "5.493"; E = "  +1.91%  " }
    15 = @{ D = "94.99"; E = "  +1.07%  " }
    16 = @{ D = "15.19"; E = "  +4.82%  " }
    17 = @{ D = "30.342.06"; E = "  +0.25%  " }
    18 = @{ D = "0.000008258"; E = "  +5.68%  " }
    19 = @{ D = "253.65"; E = "  -2.26%  " }
    20 = @{ D = "5.838"; E = "  -0.07%  " }
    21 = @{ D = "2.188.42"; E = "  +0.49%  " }
    22 = @{ D = "0.9997"; E = "  -0.16%  " }
    23 = @{ D = "0.9997"; E = "  -0.23%  " }
    24 = @{ D = "6.975"; E = "  +2.41%  " }
    25 = @{ D = "9.784"; E = "  +1.87%  " }
    26 = @{ D = "165.81"; E = "  +4.54%  " }
    27 = @{ D = "2.355"; E = "  +5.97%  " }
    28 = @{ D = "19.38" }
    29 = @{ D = "0.1306"; E = "  -0.56%  " }
    30 = @{ D = "1.353"; E = "  +0.98%  " }
    31 = @{ D = "1.541"; E = "  -0.56%  " }
    32 = @{ D = "4.449"; E = "  +1.27%  " }
    33 = @{ D = "4.221"; E = "  +1.04%  " }
    34 = @{ D = "0.05266"; E = "  +4.59%  " }
    35 = @{ E = "  +6.76%  " }
    36 = @{ D = "0.7521"; E = "  +0.85%  " }
    37 = @{ D = "2.768"; E = "  +2.24%  " }
    38 = @{ D = "0.01978"; E = "  +3.14%  " }
    39 = @{ D = "2.806"; E = "  +0.26%  " }
    40 = @{ D = "79.54"; E = "  -0.58%  " }
    41 = @{ D = "6.457"; E = "  -0.83%  " }
    42 = @{ D = "0.4567"; E = "  +3.85%  " }
    43 = @{ E = "  -0.37%  " }
    44 = @{ D = "0.8452"; E = "  +1.45%  " }
    45 = @{ D = "1.001"; E = "  -0.02%  " }
    46 = @{ D = "102.02"; E = "  +0.44%  " }
    47 = @{ D = "9.786"; E = "  +0.73%  " }
    48 = @{ D = "7.466"; E = "  +3.72%  " }
    49 = @{ D = "36.77"; E = "  +2.81%  " }
    50 = @{ D = "0.4218"; E = "  +4.16%  " }
    51 = @{ B = "Cronos"; C = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; D = "0.06049"; E = "  +1.81%  " }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $cell = $ws.Range("$col$row")
        if ($col -eq "D") {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $cols[$col]
    }
}

$wb.Save()
